$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextCell "D2" "28.273.44"
$ws.Range("E2").Value = "  -0.31%  "
Set-TextCell "D3" "1.809.20"
$ws.Range("E3").Value = "  -0.56%  "
Set-TextCell "D4" "1.002"
$ws.Range("E4").Value = "  -0.21%  "
Set-TextCell "D5" "312.37"
$ws.Range("E5").Value = "  -0.96%  "
Set-TextCell "D6" "1.001"
$ws.Range("E6").Value = "  -0.25%  "
Set-TextCell "D7" "0.5127"
$ws.Range("E7").Value = "  -2.92%  "
Set-TextCell "D8" "0.3934"
$ws.Range("E8").Value = "  +2.21%  "
Set-TextCell "D9" "0.07811"
$ws.Range("E9").Value = "  -2.76%  "
Set-TextCell "D10" "1.108"
$ws.Range("E10").Value = "  -0.28%  "
Set-TextCell "D11" "41.05"
$ws.Range("E11").Value = "  -1.93%  "
Set-TextCell "D12" "6.362"
$ws.Range("E12").Value = "  -0.29%  "
Set-TextCell "D13" "1.002"
$ws.Range("E13").Value = "  -0.16%  "
Set-TextCell "D14" "20.49"
$ws.Range("E14").Value = "  -1.81%  "
Set-TextCell "D15" "7.334"
$ws.Range("E15").Value = "  -1.06%  "
Set-TextCell "D16" "1.798.81"
$ws.Range("E16").Value = "  -1.13%  "
Set-TextCell "D17" "92.58"
$ws.Range("E17").Value = "  -1.76%  "
Set-TextCell "D18" "0.00001078"
$ws.Range("E18").Value = "  -1.98%  "
Set-TextCell "D19" "0.06563"
$ws.Range("E19").Value = "  -1.09%  "
Set-TextCell "D20" "1.000"
$ws.Range("E20").Value = "  -0.25%  "
Set-TextCell "D21" "17.31"
$ws.Range("E21").Value = "  -1.71%  "
Set-TextCell "D22" "6.015"
$ws.Range("E22").Value = "  +0.09%  "
Set-TextCell "D23" "28.321.05"
$ws.Range("E23").Value = "  -0.32%  "
Set-TextCell "D24" "11.14"
$ws.Range("E24").Value = "  -1.64%  "
Set-TextCell "D25" "2.232"
$ws.Range("E25").Value = "  -0.57%  "
Set-TextCell "D26" "161.04"
$ws.Range("E26").Value = "  +1.21%  "
Set-TextCell "D27" "2.454"
$ws.Range("E27").Value = "  +2.23%  "
Set-TextCell "D28" "20.49"
$ws.Range("E28").Value = "  -1.57%  "
Set-TextCell "D29" "2.015.15"
$ws.Range("E29").Value = "  -0.58%  "
Set-TextCell "D30" "127.67"
$ws.Range("E30").Value = "  +2.68%  "
Set-TextCell "D31" "0.1093"
$ws.Range("E31").Value = "  -1.45%  "
Set-TextCell "D32" "1.062"
$ws.Range("E32").Value = "  -1.34%  "
Set-TextCell "D33" "3.653"
$ws.Range("E33").Value = "  -0.66%  "
Set-TextCell "D34" "5.567"
$ws.Range("E34").Value = "  -1.75%  "
Set-TextCell "D35" "0.07158"
$ws.Range("E35").Value = "  -2.39%  "
Set-TextCell "D36" "9.169"
$ws.Range("E36").Value = "  +5.25%  "
Set-TextCell "D37" "0.02354"
$ws.Range("E37").Value = "  +0.65%  "
Set-TextCell "D38" "0.2178"
$ws.Range("E38").Value = "  -0.71%  "
Set-TextCell "D39" "11.56"
$ws.Range("E39").Value = "  -5.29%  "
Set-TextCell "D40" "5.017"
$ws.Range("E40").Value = "  -2.21%  "
Set-TextCell "D41" "0.6178"
$ws.Range("E41").Value = "  -1.90%  "
Set-TextCell "D42" "1.0000"
$ws.Range("E42").Value = "  -0.24%  "
Set-TextCell "D43" "1.155"
$ws.Range("E43").Value = "  -2.23%  "
Set-TextCell "D44" "13.23"
$ws.Range("E44").Value = "  -0.99%  "
Set-TextCell "D45" "0.5960"
$ws.Range("E45").Value = "  -2.41%  "
$ws.Range("E46").Value = "  -5.76%  "
Set-TextCell "D47" "3.731"
$ws.Range("E47").Value = "  -1.39%  "
Set-TextCell "D48" "125.24"
$ws.Range("E48").Value = "  -1.41%  "
Set-TextCell "D49" "1.210"
$ws.Range("E49").Value = "  +0.38%  "
Set-TextCell "D50" "1.924"
$ws.Range("E50").Value = "  -2.25%  "
Set-TextCell "D51" "0.06803"
$ws.Range("E51").Value = "  -1.33%  "
